# The workbook's "ENAME" sheet had a running list in column A (A1:A13,
# with A13 = SUM(A4:A12)). The edit selects A2:A10 and clears their
# contents (Delete key), leaving only A1, A11, A12 and the A13 total,
# which recalculates automatically (787 + 677 = 1464).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENAME")

# Clear the contents of A2:A10 (values, the SUM(A1:A2) formula in A3,
# and the shared-string text in A7 all go away).
$rng = $ws.Range("A2:A10")
$rng.ClearContents()

# Leave the selection where the user left it after the delete: A2:A10
# with the active cell at A2.
[void]$rng.Select()
